# Insert a new price record as row 158 in the daily "Camote" price log.
# Excel shifts the existing rows 158:227 down to 159:228 and we populate
# the newly opened row 158 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(158).EntireRow.Insert()

$ws.Range("A158").Value = 10
$ws.Range("B158").Value = "Vega Modelo de Temuco"
$ws.Range("C158").Value = "La Araucanía"
$ws.Range("D158").Value = 45141
$ws.Range("E158").Value = 9
$ws.Range("F158").Value = 100114002
$ws.Range("G158").Value = "Camote"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 100
$ws.Range("K158").Value = 26000
$ws.Range("L158").Value = 26000
$ws.Range("M158").Value = 26000
$ws.Range("N158").Value = "$/caja 18 kilos"
$ws.Range("O158").Value = "Perú"
$ws.Range("P158").Value = 1444
$ws.Range("Q158").Value = 18
$ws.Range("R158").Value = "Hortaliza"
